# feat: add 2022-Q1 data
#
# The workbook currently has two sheets: "2021-Q4" (fund holdings detail)
# and "总计" (per-quarter summary). We need to:
#   1. Insert a new "2022-Q1" fund-holdings sheet between them, populated
#      with the latest quarter's fund data.
#   2. Update the "总计" summary sheet with a new leading row for 2022-Q1
#      (the old "总计" sheet slot gets reused/renamed as "2022-Q1", and a
#      fresh "总计" sheet is appended after it, matching how the workbook's
#      sheetId/rId numbering comes out after the edit).

$wb = $excel.ActiveWorkbook
$detail = $wb.Worksheets.Item(1)   # "2021-Q4" - used as a formatting template
$total = $wb.Worksheets.Item(2)    # currently "总计"

# ---------------------------------------------------------------------
# Step 1: turn the existing "总计" sheet into the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$total.Cells.Clear()
$total.Name = "2022-Q1"

# Copy header-row / first-column formatting from the "2021-Q4" sheet so the
# new sheet matches the workbook's existing look (bold, centered, bordered).
$detail.Range("B1:H1").Copy()
$total.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$detail.Range("A2").Copy()
$total.Range("A2").PasteSpecial(-4122)      # xlPasteFormats

$total.Range("A2").Value = 0

$total.Range("B1").Value = "基金代码"
$total.Range("C1").Value = "基金名称"
$total.Range("D1").Value = "基金规模"
$total.Range("E1").Value = "股票总仓位"
$total.Range("F1").Value = "仓位占比"
$total.Range("G1").Value = "持有市值(亿元)"
$total.Range("H1").Value = "仓位排名"

# B2:G2 hold text-like values (fund code / numbers-as-text), so force text
# formatting before assigning, then drop back to the Normal style so no
# extra per-cell formatting is left behind.
$fundDataRange = $total.Range("B2:G2")
$fundDataRange.NumberFormat = "@"
$total.Range("B2").Value = "161123"
$total.Range("C2").Value = "易方达并购重组指数（LOF）"
$total.Range("D2").Value = "4.78"
$total.Range("E2").Value = "94.71"
$total.Range("F2").Value = "3.62"
$total.Range("G2").Value = "0.1730"
$fundDataRange.Style = "Normal"

$total.Range("H2").Value = 8

# ---------------------------------------------------------------------
# Step 2: append a brand-new "总计" sheet after "2022-Q1" with the
# updated summary (2022-Q1 on top, 2021-Q4 below).
# ---------------------------------------------------------------------
$detail.Copy($null, $total)
$grand = $wb.Worksheets.Item(3)
$grand.Cells.Clear()
$grand.Name = "总计"

$detail.Range("B1:D1").Copy()
$grand.Range("B1:D1").PasteSpecial(-4122)   # xlPasteFormats
$detail.Range("A2:A3").Copy()
$grand.Range("A2:A3").PasteSpecial(-4122)   # xlPasteFormats

$grand.Range("B1").Value = "日期"
$grand.Range("C1").Value = "持有数量(只)"
$grand.Range("D1").Value = "持有市值(亿元)"

$grand.Range("A2").Value = 0
$grand.Range("B2").Value = "2022-Q1"
$grand.Range("C2").Value = 1
$grand.Range("D2").Value = 0.17

$grand.Range("A3").Value = 1
$grand.Range("B3").Value = "2021-Q4"
$grand.Range("C3").Value = 3
$grand.Range("D3").Value = 0.3

# Restore the originally active sheet/selection so we don't leave an
# unrelated UI-state change behind.
$detail.Activate()
$detail.Range("A1").Select() | Out-Null
